$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric-looking Price/Volume cells so they
# keep their exact literal representation instead of being parsed as numbers.
$numCells = @("D2","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","E23","D24","E24","D25","E25","D26","E26","D27","E27","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47")
foreach ($addr in $numCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) / Volume(1h) (E) columns
$ws.Range("D2").Value = "275.29"
$ws.Range("E3").Value = "-2.20%"
$ws.Range("D4").Value = "4.858"
$ws.Range("E4").Value = "1.26%"
$ws.Range("D5").Value = "0.06349"
$ws.Range("E5").Value = "1.32%"
$ws.Range("D6").Value = "6.884"
$ws.Range("E6").Value = "-0.58%"
$ws.Range("D7").Value = "3.313"
$ws.Range("E7").Value = "1.43%"
$ws.Range("D8").Value = "1.257"
$ws.Range("E8").Value = "33.30%"
$ws.Range("D9").Value = "0.8687"
$ws.Range("E9").Value = "-1.14%"
$ws.Range("D10").Value = "0.1576"
$ws.Range("E10").Value = "8.07%"
$ws.Range("D11").Value = "0.05197"
$ws.Range("E11").Value = "-2.60%"
$ws.Range("D12").Value = "0.07393"
$ws.Range("E12").Value = "1.36%"
$ws.Range("D13").Value = "0.02936"
$ws.Range("E13").Value = "-5.62%"
$ws.Range("D14").Value = "0.09030"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("E15").Value = "2.11%"
$ws.Range("D16").Value = "0.0006333"
$ws.Range("E16").Value = "1.12%"
$ws.Range("D17").Value = "0.005943"
$ws.Range("E17").Value = "2.50%"
$ws.Range("D18").Value = "3.449"
$ws.Range("E18").Value = "0.13%"
$ws.Range("D19").Value = "2.271"
$ws.Range("E19").Value = "-0.55%"
$ws.Range("D20").Value = "0.3114"
$ws.Range("E20").Value = "-1.06%"
$ws.Range("D21").Value = "0.1334"
$ws.Range("E21").Value = "1.67%"
$ws.Range("D22").Value = "3.903"
$ws.Range("E22").Value = "1.41%"
$ws.Range("E23").Value = "0.72%"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").Value = "0.12%"
$ws.Range("D25").Value = "0.004251"
$ws.Range("E25").Value = "-0.61%"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").Value = "-0.21%"
$ws.Range("D27").Value = "0.0001678"
$ws.Range("E27").Value = "-0.77%"
$ws.Range("D40").Value = "0.04119"
$ws.Range("E40").Value = "1.84%"
$ws.Range("D41").Value = "0.006796"
$ws.Range("E41").Value = "6.12%"
$ws.Range("D42").Value = "0.1166"
$ws.Range("E42").Value = "1.10%"
$ws.Range("D43").Value = "0.002144"
$ws.Range("E43").Value = "-0.18%"
$ws.Range("D44").Value = "0.01074"
$ws.Range("E44").Value = "-10.15%"
$ws.Range("D45").Value = "0.00005304"
$ws.Range("E45").Value = "4.35%"
$ws.Range("D46").Value = "1.486"
$ws.Range("E46").Value = "-37.47%"
$ws.Range("D47").Value = "0.02099"
$ws.Range("E47").Value = "-29.72%"

# Update Coin (B) / Link (C) columns (rows 46 and 47 swap content)
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
